$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.360.95'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '2.912.22'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '361.95'
$c.ClearFormats()
$ws.Range('E5').Value = '  +2.14%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '104.72'
$c.ClearFormats()
$ws.Range('E6').Value = '  -3.65%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.544'
$c.ClearFormats()
$ws.Range('E7').Value = '  -3.67%  '
$ws.Range('E8').Value = '  -0.10%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.590'
$c.ClearFormats()
$ws.Range('E9').Value = '  -5.11%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.92'
$c.ClearFormats()
$ws.Range('E10').Value = '  -4.65%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.140'
$c.ClearFormats()
$ws.Range('E11').Value = '  +1.76%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0836'
$c.ClearFormats()
$ws.Range('E12').Value = '  -3.40%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '18.57'
$c.ClearFormats()
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('D14').Value = '3.376.37'
$ws.Range('E14').Value = '  +0.38%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.35'
$c.ClearFormats()
$ws.Range('E15').Value = '  -4.29%  '
$ws.Range('D16').Value = '2.930.59'
$ws.Range('E16').Value = '  +0.73%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.961'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '51.278.33'
$ws.Range('E18').Value = '  -0.73%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '3.29'
$c.ClearFormats()
$ws.Range('E19').Value = '  -2.24%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.23'
$c.ClearFormats()
$ws.Range('E20').Value = '  -3.49%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.04'
$c.ClearFormats()
$ws.Range('E21').Value = '  -5.48%  '
$ws.Range('D22').Value = '0.0₃0946'
$ws.Range('E22').Value = '  -2.70%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '68.40'
$c.ClearFormats()
$ws.Range('E23').Value = '  -2.50%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '260.39'
$c.ClearFormats()
$ws.Range('E24').Value = '  -2.51%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.70'
$c.ClearFormats()
$ws.Range('E25').Value = '  -3.43%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.175'
$c.ClearFormats()
$ws.Range('E26').Value = '  -4.45%  '
$ws.Range('E27').Value = '  +0.01%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '26.12'
$c.ClearFormats()
$ws.Range('E28').Value = '  -2.14%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.28'
$c.ClearFormats()
$ws.Range('E29').Value = '  -3.27%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.109'
$c.ClearFormats()
$ws.Range('E30').Value = '  +3.94%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '10.00'
$c.ClearFormats()
$ws.Range('E31').Value = '  -4.25%  '
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '6.14'
$c.ClearFormats()
$ws.Range('E32').Value = '  +1.33%  '
$ws.Range('E33').Value = '  -2.34%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '35.00'
$c.ClearFormats()
$ws.Range('E34').Value = '  -5.65%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '51.08'
$c.ClearFormats()
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('E36').Value = '  +0.41%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0424'
$c.ClearFormats()
$ws.Range('E37').Value = '  -3.23%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.81'
$c.ClearFormats()
$ws.Range('E38').Value = '  +4.88%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.16'
$c.ClearFormats()
$ws.Range('E39').Value = '  -0.29%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '17.04'
$c.ClearFormats()
$ws.Range('E40').Value = '  -5.89%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.87'
$c.ClearFormats()
$ws.Range('E41').Value = '  -5.84%  '
$ws.Range('E42').Value = '  -3.66%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '22.47'
$c.ClearFormats()
$ws.Range('E43').Value = '  -1.46%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '119.77'
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.75%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.14'
$c.ClearFormats()
$ws.Range('E45').Value = '  -1.40%  '
$ws.Range('D46').Value = '2.080.16'
$ws.Range('E46').Value = '  -1.84%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.20'
$c.ClearFormats()
$ws.Range('E47').Value = '  -6.30%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.ClearFormats()
$ws.Range('E48').Value = '  -8.50%  '
$ws.Range('D49').Value = '3.212.44'
$ws.Range('E49').Value = '  +0.55%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.238'
$c.ClearFormats()
$ws.Range('E50').Value = '  -3.97%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0313'
$c.ClearFormats()
$ws.Range('E51').Value = '  -7.66%  '
